$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-06-12 Thursday" "2025-06-13 Friday"

Replace-Text "304÷9=" "477÷8="
Replace-Text "173÷7=" "629÷7="
Replace-Text "779÷5=" "491÷5="
Replace-Text "382÷6=" "826÷8="
Replace-Text "882÷4=" "855÷8="

Replace-Text "947÷5=" "775÷3="
Replace-Text "391÷4=" "454÷9="
Replace-Text "809÷7=" "842÷8="
Replace-Text "530÷5=" "865÷7="
Replace-Text "558÷6=" "664÷9="

Replace-Text "106÷4=" "244÷3="
Replace-Text "643÷7=" "717÷7="
Replace-Text "805÷6=" "699÷6="
Replace-Text "732÷5=" "918÷9="
Replace-Text "876÷6=" "450÷9="

Replace-Text "421÷6=" "766÷3="
Replace-Text "414÷8=" "916÷4="
Replace-Text "681÷3=" "986÷7="
Replace-Text "649÷8=" "136÷5="
Replace-Text "127÷2=" "206÷9="

Replace-Text "855÷9=" "628÷9="
Replace-Text "746÷4=" "713÷7="
Replace-Text "658÷2=" "296÷7="
Replace-Text "416÷2=" "989÷8="
Replace-Text "453÷5=" "882÷9="
